$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.381467440838577
$ws.Range("C2").Value = 0.08471188263121121
$ws.Range("D2").Value = 0.07817739636989529
$ws.Range("E2").Value = 0.1321307422917499
$ws.Range("G2").Value = 1.063198050451817
$ws.Range("H2").Value = 1.039708937878558
$ws.Range("I2").Value = 0.84507816982984
$ws.Range("K2").Value = 0.4499417794007172
$ws.Range("M2").Value = 0.2696813439553907
$ws.Range("N2").Value = 1.907622604444569

$ws.Range("B3").Value = 0.3451514959369888
$ws.Range("C3").Value = 0.07535324325306192
$ws.Range("D3").Value = 0.07095966224380845
$ws.Range("E3").Value = 0.1220732104643147
$ws.Range("G3").Value = 1.050880262502631
$ws.Range("H3").Value = 1.038839395086342
$ws.Range("I3").Value = 0.8435318156932468
$ws.Range("K3").Value = 0.4052718456768787
$ws.Range("M3").Value = 0.2458305003562629
$ws.Range("N3").Value = 1.922175776200191

$ws.Range("B4").Value = 0.3229969342560253
$ws.Range("C4").Value = 0.06963414685066027
$ws.Range("D4").Value = 0.06656436765099727
$ws.Range("E4").Value = 0.1159804873851797
$ws.Range("G4").Value = 1.043956324700432
$ws.Range("H4").Value = 1.038777973989994
$ws.Range("I4").Value = 0.8430047865315444
$ws.Range("K4").Value = 0.3780086838154659
$ws.Range("M4").Value = 0.2313205919215875
$ws.Range("N4").Value = 1.93172775197548

$ws.Range("B5").Value = 0.314004944827758
$ws.Range("C5").Value = 0.06731030431990348
$ws.Range("D5").Value = 0.06478236409206772
$ws.Range("E5").Value = 0.1135181912612993
$ws.Range("G5").Value = 1.041295129104142
$ws.Range("H5").Value = 1.038871695390867
$ws.Range("I5").Value = 0.8428961102187316
$ws.Range("K5").Value = 0.3669400082538914
$ws.Range("M5").Value = 0.2254412709067637
$ws.Range("N5").Value = 1.935775044775248

$ws.Range("B6").Value = 0.3125140188705302
$ws.Range("C6").Value = 0.06692483623970702
$ws.Range("D6").Value = 0.06448701298201343
$ws.Range("E6").Value = 0.1131105622322082
$ws.Range("G6").Value = 1.040862914366102
$ws.Range("H6").Value = 1.038894428370327
$ws.Range("I6").Value = 0.8428844684582657
$ws.Range("K6").Value = 0.3651045553943959
$ws.Range("M6").Value = 0.2244670349157616
$ws.Range("N6").Value = 1.936456438059899

$ws.Range("B7").Value = 0.3228755185766659
$ws.Range("C7").Value = 0.06960277954374305
$ws.Range("D7").Value = 0.066540298080497
$ws.Range("E7").Value = 0.1159471971605441
$ws.Range("G7").Value = 1.043919786073261
$ws.Range("H7").Value = 1.038778757213322
$ws.Range("I7").Value = 0.8430028915008378
$ws.Range("K7").Value = 0.3778592406837902
$ws.Range("M7").Value = 0.2312411657052778
$ws.Range("N7").Value = 1.931781708710851

$ws.Range("B8").Value = 0.3689159342582684
$ws.Range("C8").Value = 0.08147933219595416
$ws.Range("D8").Value = 0.07568113300771984
$ws.Range("E8").Value = 0.1286456032605159
$ws.Range("G8").Value = 1.058817971131347
$ws.Range("H8").Value = 1.039311017553118
$ws.Range("I8").Value = 0.8444572266712456
$ws.Range("K8").Value = 0.4345054315976995
$ws.Range("M8").Value = 0.2614294691759866
$ws.Range("N8").Value = 1.91251252290045

$ws.Range("B9").Value = 0.4603421249071857
$ws.Range("C9").Value = 0.1049898898190236
$ws.Range("D9").Value = 0.09389824754859433
$ws.Range("E9").Value = 0.1542151704746146
$ws.Range("G9").Value = 1.093125202984012
$ws.Range("H9").Value = 1.0441076876522
$ws.Range("I9").Value = 0.8506682682202182
$ws.Range("K9").Value = 0.54690068364377
$ws.Range("M9").Value = 0.32171176364281
$ws.Range("N9").Value = 1.879623441515427

$ws.Range("B10").Value = 0.5282184908800787
$ws.Range("C10").Value = 0.1224066745021162
$ws.Range("D10").Value = 0.1074658426261692
$ws.Range("E10").Value = 0.1734273458002846
$ws.Range("G10").Value = 1.121466752371333
$ws.Range("H10").Value = 1.049927131466916
$ws.Range("I10").Value = 0.8572910982319897
$ws.Range("K10").Value = 0.6302975284988861
$ws.Range("M10").Value = 0.3666869756999134
$ws.Range("N10").Value = 1.858456226358975

$ws.Range("B11").Value = 0.559253013118564
$ws.Range("C11").Value = 0.1303632569571107
$ws.Range("D11").Value = 0.113679198544034
$ws.Range("E11").Value = 0.1822642252568158
$ws.Range("G11").Value = 1.135048181948378
$ws.Range("H11").Value = 1.053074728582089
$ws.Range("I11").Value = 0.8607538940568062
$ws.Range("K11").Value = 0.6684198625530371
$ws.Range("M11").Value = 0.387302071746646
$ws.Range("N11").Value = 1.849479695325478

$ws.Range("B12").Value = 0.5710276398963288
$ws.Range("C12").Value = 0.133381168247837
$ws.Range("D12").Value = 0.1160380574261239
$ws.Range("E12").Value = 0.1856247948169241
$ws.Range("G12").Value = 1.140290652053608
$ws.Range("H12").Value = 1.054338698570007
$ws.Range("I12").Value = 0.8621300682116271
$ws.Range("K12").Value = 0.6828825606411328
$ws.Range("M12").Value = 0.3951312230455599
$ws.Range("N12").Value = 1.8461745711597

$ws.Range("B13").Value = 0.5684907627394864
$ws.Range("C13").Value = 0.1327309864442157
$ws.Range("D13").Value = 0.1155297675892797
$ws.Range("E13").Value = 0.1849003987188098
$ws.Range("G13").Value = 1.139157160973127
$ws.Range("H13").Value = 1.054063274249899
$ws.Range("I13").Value = 0.8618307957620743
$ws.Range("K13").Value = 0.6797665762503016
$ws.Range("M13").Value = 0.3934440628576255
$ws.Range("N13").Value = 1.846882199794685

$ws.Range("B14").Value = 0.5602212677411558
$ws.Range("C14").Value = 0.1306114430618948
$ws.Range("D14").Value = 0.1138731428643354
$ws.Range("E14").Value = 0.182540414532042
$ws.Range("G14").Value = 1.1354774869057
$ws.Range("H14").Value = 1.053177271716692
$ws.Range("I14").Value = 0.8608658112683472
$ws.Range("K14").Value = 0.6696091841401994
$ws.Range("M14").Value = 0.3879457249466327
$ws.Range("N14").Value = 1.849205893505456

$ws.Range("B15").Value = 0.5551588945352819
$ws.Range("C15").Value = 0.129313805898164
$ws.Range("D15").Value = 0.1128591945358579
$ws.Range("E15").Value = 0.1810967177381855
$ws.Range("G15").Value = 1.133236548334139
$ws.Range("H15").Value = 1.052643954759475
$ws.Range("I15").Value = 0.8602831859921167
$ws.Range("K15").Value = 0.6633909531797428
$ws.Range("M15").Value = 0.3845807938050143
$ws.Range("N15").Value = 1.850641484471005

$ws.Range("B16").Value = 0.5261934685585743
$ws.Range("C16").Value = 0.1218873773064217
$ws.Range("D16").Value = 0.1070606218274435
$ws.Range("E16").Value = 0.1728518125832181
$ws.Range("G16").Value = 1.120593070360826
$ws.Range("H16").Value = 1.049731504489301
$ws.Range("I16").Value = 0.8570738670878413
$ws.Range("K16").Value = 0.6278098640577241
$ws.Range("M16").Value = 0.3653428865858075
$ws.Range("N16").Value = 1.859056016592277

$ws.Range("B17").Value = 0.5084643276563554
$ws.Range("C17").Value = 0.1173401766546931
$ws.Range("D17").Value = 0.1035140173672886
$ws.Range("E17").Value = 0.1678189030439938
$ws.Range("G17").Value = 1.113013415682843
$ws.Range("H17").Value = 1.048073020818435
$ws.Range("I17").Value = 0.8552204422421141
$ws.Range("K17").Value = 0.6060293491881055
$ws.Range("M17").Value = 0.3535811169850049
$ws.Range("N17").Value = 1.864385373781417

$ws.Range("B18").Value = 0.4982817980488221
$ws.Range("C18").Value = 0.1147279039739431
$ws.Range("D18").Value = 0.101477997358856
$ws.Range("E18").Value = 0.1649332363617617
$ws.Range("G18").Value = 1.108718597485336
$ws.Range("H18").Value = 1.047166189592446
$ws.Range("I18").Value = 0.8541967504090593
$ws.Range("K18").Value = 0.5935191545970326
$ws.Range("M18").Value = 0.3468306915150166
$ws.Range("N18").Value = 1.867512083602861

$ws.Range("B19").Value = 0.4948367134921909
$ws.Range("C19").Value = 0.113843972654422
$ws.Range("D19").Value = 0.1007893032610667
$ws.Range("E19").Value = 0.1639577592661396
$ws.Range("G19").Value = 1.107275561702522
$ws.Range("H19").Value = 1.046867235954608
$ws.Range("I19").Value = 0.8538574141883046
$ws.Range("K19").Value = 0.5892864017990007
$ws.Range("M19").Value = 0.3445476182257039
$ws.Range("N19").Value = 1.868581272306649

$ws.Range("B20").Value = 0.5103500934833392
$ws.Range("C20").Value = 0.1178239065213234
$ws.Range("D20").Value = 0.1038911562290394
$ws.Range("E20").Value = 0.1683537180624839
$ws.Range("G20").Value = 1.113813573252912
$ws.Range("H20").Value = 1.048244695419811
$ws.Range("I20").Value = 0.8554133582495353
$ws.Range("K20").Value = 0.6083461229868021
$ws.Range("M20").Value = 0.3548316609261519
$ws.Range("N20").Value = 1.863811697601271

$ws.Range("B21").Value = 0.562649608148746
$ws.Range("C21").Value = 0.1312338698793383
$ws.Range("D21").Value = 0.1143595706257088
$ws.Range("E21").Value = 0.1832332106540235
$ws.Range("G21").Value = 1.136555593055988
$ws.Range("H21").Value = 1.053435556092523
$ws.Range("I21").Value = 0.8611474881556305
$ws.Range("K21").Value = 0.6725919342460429
$ws.Range("M21").Value = 0.3895601027941851
$ws.Range("N21").Value = 1.84852081258795

$ws.Range("B22").Value = 0.59696163101097
$ws.Range("C22").Value = 0.1400268142569985
$ws.Range("D22").Value = 0.1212362535360114
$ws.Range("E22").Value = 0.193040915529842
$ws.Range("G22").Value = 1.151998855455105
$ws.Range("H22").Value = 1.057248055176302
$ws.Range("I22").Value = 0.8652733755061917
$ws.Range("K22").Value = 0.714735336835048
$ws.Range("M22").Value = 0.4123893692304819
$ws.Range("N22").Value = 1.839075906082343

$ws.Range("B23").Value = 0.5786366741901929
$ws.Range("C23").Value = 0.1353311915183895
$ws.Range("D23").Value = 0.1175628241602027
$ws.Range("E23").Value = 0.1877986681960664
$ws.Range("G23").Value = 1.143703275673033
$ws.Range("H23").Value = 1.055174791489947
$ws.Range("I23").Value = 0.8630366392355384
$ws.Range("K23").Value = 0.6922284088849153
$ws.Range("M23").Value = 0.4001927686270648
$ws.Range("N23").Value = 1.84406654967664

$ws.Range("B24").Value = 0.5094975070096837
$ws.Range("C24").Value = 0.1176052061101132
$ws.Range("D24").Value = 0.1037206425038164
$ws.Range("E24").Value = 0.1681119038582466
$ws.Range("G24").Value = 1.113451626342567
$ws.Range("H24").Value = 1.048166936008101
$ws.Range("I24").Value = 0.855326010515931
$ws.Range("K24").Value = 0.6072986730733874
$ws.Range("M24").Value = 0.3542662539904171
$ws.Range("N24").Value = 1.864070861011427

$ws.Range("B25").Value = 0.4354857413298703
$ws.Range("C25").Value = 0.09860510391322919
$ws.Range("D25").Value = 0.0889382113253987
$ws.Range("E25").Value = 0.1472244532654017
$ws.Range("G25").Value = 1.083295781067108
$ws.Range("H25").Value = 1.042407524555159
$ws.Range("I25").Value = 0.8486271917473616
$ws.Range("K25").Value = 0.5163522404556886
$ws.Range("M25").Value = 0.3052852209564278
$ws.Range("N25").Value = 1.887995392565131
